$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.790.66'
$ws.Range('E2').Value = '  -1.77%  '
$ws.Range('D3').Value = '3.133.20'
$ws.Range('E3').Value = '  -7.35%  '
$ws.Range('E4').Value = '  +0.01%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '565.90'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  -2.52%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '167.77'
$cell.Style = "Normal"
$ws.Range('E6').Value = '  -6.15%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '0.607'
$cell.Style = "Normal"
$ws.Range('E7').Value = '  -1.81%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '3.135.05'
$ws.Range('E9').Value = '  -7.25%  '
$ws.Range('E10').Value = '  -5.60%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '6.51'
$cell.Style = "Normal"
$ws.Range('E11').Value = '  -5.82%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '0.388'
$cell.Style = "Normal"
$ws.Range('E12').Value = '  -5.08%  '
$ws.Range('D13').Value = '3.680.88'
$ws.Range('E13').Value = '  -7.57%  '
$ws.Range('E14').Value = '  +1.12%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '26.70'
$cell.Style = "Normal"
$ws.Range('E15').Value = '  -7.57%  '
$ws.Range('D16').Value = '64.719.56'
$ws.Range('E16').Value = '  -2.06%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '0.0000160'
$cell.Style = "Normal"
$ws.Range('E17').Value = '  -6.24%  '
$ws.Range('D18').Value = '3.142.64'
$ws.Range('E18').Value = '  -7.65%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '5.67'
$cell.Style = "Normal"
$ws.Range('E19').Value = '  -3.06%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '12.66'
$cell.Style = "Normal"
$ws.Range('E20').Value = '  -7.55%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '354.55'
$cell.Style = "Normal"
$ws.Range('E21').Value = '  -2.82%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '7.17'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  -4.69%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  +0.37%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '68.81'
$cell.Style = "Normal"
$ws.Range('E24').Value = '  -5.28%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '3.291.04'
$ws.Range('E25').Value = '  -7.49%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '0.491'
$cell.Style = "Normal"
$ws.Range('E26').Value = '  -6.99%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '0.0000114'
$cell.Style = "Normal"
$ws.Range('E27').Value = '  -6.92%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '9.57'
$cell.Style = "Normal"
$ws.Range('E28').Value = '  -1.87%  '
$ws.Range('E29').Value = '  -1.57%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('E32').Value = '  -4.57%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '21.74'
$cell.Style = "Normal"
$ws.Range('E33').Value = '  -5.79%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '5.25'
$cell.Style = "Normal"
$ws.Range('E34').Value = '  -8.27%  '
$ws.Range('E35').Value = '  -4.60%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '6.53'
$cell.Style = "Normal"
$ws.Range('E36').Value = '  -6.40%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '157.58'
$cell.Style = "Normal"
$ws.Range('E37').Value = '  -3.05%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '1.42'
$cell.Style = "Normal"
$ws.Range('E38').Value = '  -6.99%  '
$ws.Range('E39').Value = '  -3.11%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '1.75'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  -1.76%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '25.75'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  -5.25%  '
$ws.Range('D42').Value = '2.647.78'
$ws.Range('E42').Value = '  -1.25%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '2.41'
$cell.Style = "Normal"
$ws.Range('E43').Value = '  -7.02%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '6.05'
$cell.Style = "Normal"
$ws.Range('E44').Value = '  -2.60%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '4.12'
$cell.Style = "Normal"
$ws.Range('E45').Value = '  -4.86%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '39.34'
$cell.Style = "Normal"
$ws.Range('E46').Value = '  -0.67%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '0.0648'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  -4.74%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '23.73'
$cell.Style = "Normal"
$ws.Range('E48').Value = '  -3.41%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '316.41'
$cell.Style = "Normal"
$ws.Range('E49').Value = '  -4.31%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '0.0269'
$cell.Style = "Normal"
$ws.Range('E50').Value = '  -5.11%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '0.102'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  -1.49%  '
